$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Split the Huckemann et al. reference run so that the journal name
# "Bernoulli" is italicized, matching:
#   " ... differentiation. " + italic("Bernoulli") + ", 22(4), 2113-2142."
# ---------------------------------------------------------------------------
$target = $d.Content
$found = $target.Find.Execute("Bernoulli, 22(4), 2113-2142.")
if (-not $found) {
    throw "Could not find the Bernoulli citation text to split."
}
$citeStart = $target.Start
$citeEnd = $target.End

$splitXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Bernoulli</w:t></w:r><w:r><w:t>, 22(4), 2113-2142.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($splitXml) | Out-Null

# InsertXML inserts (rather than replaces), so the original, non-italicized
# "Bernoulli, 22(4), 2113-2142." text is still sitting just before the newly
# inserted (properly formatted) copy. Remove that now-duplicate plain text.
$dupRange = $d.Range($citeStart, $citeStart + 28)
if ($dupRange.Text -ne "Bernoulli, 22(4), 2113-2142.") {
    throw "Unexpected text found while removing duplicate citation text: $($dupRange.Text)"
}
$dupRange.Text = ""

# ---------------------------------------------------------------------------
# Step 2: Insert a brand-new paragraph directly after the Huckemann et al.
# paragraph for the Matuk, Kurtek, & Bharath (2021) reference, with the
# arXiv identifier italicized.
# ---------------------------------------------------------------------------
$huckemannPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Huckemann*") {
        $huckemannPara = $p
        break
    }
}
if ($null -eq $huckemannPara) {
    throw "Could not find the Huckemann et al. paragraph."
}

$huckemannPara.Range.InsertParagraphAfter()

$newParaIndex = $i + 1
$newPara = $d.Paragraphs.Item($newParaIndex)

$matukXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Matuk, J., Kurtek, S., &amp; Bharath, K. (2021). Topological Data Analysis through alignment of Persistence Landscapes. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>arXiv</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> preprint arXiv:2106.15436</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($matukXml) | Out-Null
